$d = $word.ActiveDocument

$d.Content.Find.Execute("2023-06-14 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-06-15 Thursday", 2) | Out-Null
$d.Content.Find.Execute("74×14=1036", $true, $false, $false, $false, $false, $true, 1, $false, "68×63=4284", 2) | Out-Null
$d.Content.Find.Execute("29×36=1044", $true, $false, $false, $false, $false, $true, 1, $false, "65×70=4550", 2) | Out-Null
$d.Content.Find.Execute("51×33=1683", $true, $false, $false, $false, $false, $true, 1, $false, "58×41=2378", 2) | Out-Null
$d.Content.Find.Execute("30×57=1710", $true, $false, $false, $false, $false, $true, 1, $false, "96×13=1248", 2) | Out-Null
$d.Content.Find.Execute("27×26=702", $true, $false, $false, $false, $false, $true, 1, $false, "11×14=154", 2) | Out-Null
$d.Content.Find.Execute("62×100=6200", $true, $false, $false, $false, $false, $true, 1, $false, "85×46=3910", 2) | Out-Null
$d.Content.Find.Execute("85×54=4590", $true, $false, $false, $false, $false, $true, 1, $false, "38×50=1900", 2) | Out-Null
$d.Content.Find.Execute("55×81=4455", $true, $false, $false, $false, $false, $true, 1, $false, "68×59=4012", 2) | Out-Null
$d.Content.Find.Execute("14×27=378", $true, $false, $false, $false, $false, $true, 1, $false, "22×71=1562", 2) | Out-Null
$d.Content.Find.Execute("26×96=2496", $true, $false, $false, $false, $false, $true, 1, $false, "10×75=750", 2) | Out-Null
$d.Content.Find.Execute("75×48=3600", $true, $false, $false, $false, $false, $true, 1, $false, "82×26=2132", 2) | Out-Null
$d.Content.Find.Execute("35×28=980", $true, $false, $false, $false, $false, $true, 1, $false, "33×40=1320", 2) | Out-Null
$d.Content.Find.Execute("90×59=5310", $true, $false, $false, $false, $false, $true, 1, $false, "94×34=3196", 2) | Out-Null
$d.Content.Find.Execute("67×41=2747", $true, $false, $false, $false, $false, $true, 1, $false, "68×61=4148", 2) | Out-Null
$d.Content.Find.Execute("45×93=4185", $true, $false, $false, $false, $false, $true, 1, $false, "23×36=828", 2) | Out-Null
$d.Content.Find.Execute("27×66=1782", $true, $false, $false, $false, $false, $true, 1, $false, "78×72=5616", 2) | Out-Null
$d.Content.Find.Execute("13×83=1079", $true, $false, $false, $false, $false, $true, 1, $false, "63×74=4662", 2) | Out-Null
$d.Content.Find.Execute("63×84=5292", $true, $false, $false, $false, $false, $true, 1, $false, "19×94=1786", 2) | Out-Null
$d.Content.Find.Execute("76×72=5472", $true, $false, $false, $false, $false, $true, 1, $false, "93×90=8370", 2) | Out-Null
$d.Content.Find.Execute("62×76=4712", $true, $false, $false, $false, $false, $true, 1, $false, "29×96=2784", 2) | Out-Null
$d.Content.Find.Execute("98×70=6860", $true, $false, $false, $false, $false, $true, 1, $false, "18×42=756", 2) | Out-Null
$d.Content.Find.Execute("32×78=2496", $true, $false, $false, $false, $false, $true, 1, $false, "67×88=5896", 2) | Out-Null
$d.Content.Find.Execute("25×57=1425", $true, $false, $false, $false, $false, $true, 1, $false, "41×60=2460", 2) | Out-Null
$d.Content.Find.Execute("29×92=2668", $true, $false, $false, $false, $false, $true, 1, $false, "90×52=4680", 2) | Out-Null
$d.Content.Find.Execute("65×23=1495", $true, $false, $false, $false, $false, $true, 1, $false, "98×52=5096", 2) | Out-Null
$d.Content.Find.Execute("59×95=5605", $true, $false, $false, $false, $false, $true, 1, $false, "53×72=3816", 2) | Out-Null
$d.Content.Find.Execute("65×38=2470", $true, $false, $false, $false, $false, $true, 1, $false, "94×75=7050", 2) | Out-Null
$d.Content.Find.Execute("37×19=703", $true, $false, $false, $false, $false, $true, 1, $false, "75×52=3900", 2) | Out-Null
$d.Content.Find.Execute("73×33=2409", $true, $false, $false, $false, $false, $true, 1, $false, "80×99=7920", 2) | Out-Null
$d.Content.Find.Execute("80×43=3440", $true, $false, $false, $false, $false, $true, 1, $false, "38×32=1216", 2) | Out-Null
$d.Content.Find.Execute("28×10=280", $true, $false, $false, $false, $false, $true, 1, $false, "28×28=784", 2) | Out-Null
$d.Content.Find.Execute("50×84=4200", $true, $false, $false, $false, $false, $true, 1, $false, "89×77=6853", 2) | Out-Null
$d.Content.Find.Execute("31×80=2480", $true, $false, $false, $false, $false, $true, 1, $false, "77×67=5159", 2) | Out-Null
$d.Content.Find.Execute("99×63=6237", $true, $false, $false, $false, $false, $true, 1, $false, "46×27=1242", 2) | Out-Null
$d.Content.Find.Execute("63×91=5733", $true, $false, $false, $false, $false, $true, 1, $false, "64×38=2432", 2) | Out-Null
$d.Content.Find.Execute("93×79=7347", $true, $false, $false, $false, $false, $true, 1, $false, "30×67=2010", 2) | Out-Null
$d.Content.Find.Execute("89×100=8900", $true, $false, $false, $false, $false, $true, 1, $false, "70×90=6300", 2) | Out-Null
$d.Content.Find.Execute("69×17=1173", $true, $false, $false, $false, $false, $true, 1, $false, "52×58=3016", 2) | Out-Null
$d.Content.Find.Execute("87×19=1653", $true, $false, $false, $false, $false, $true, 1, $false, "80×72=5760", 2) | Out-Null
$d.Content.Find.Execute("50×38=1900", $true, $false, $false, $false, $false, $true, 1, $false, "75×44=3300", 2) | Out-Null
$d.Content.Find.Execute("91×61=5551", $true, $false, $false, $false, $false, $true, 1, $false, "23×13=299", 2) | Out-Null
$d.Content.Find.Execute("99×82=8118", $true, $false, $false, $false, $false, $true, 1, $false, "46×80=3680", 2) | Out-Null
$d.Content.Find.Execute("21×45=945", $true, $false, $false, $false, $false, $true, 1, $false, "63×10=630", 2) | Out-Null
$d.Content.Find.Execute("31×74=2294", $true, $false, $false, $false, $false, $true, 1, $false, "70×87=6090", 2) | Out-Null
$d.Content.Find.Execute("44×70=3080", $true, $false, $false, $false, $false, $true, 1, $false, "71×43=3053", 2) | Out-Null
$d.Content.Find.Execute("12×25=300", $true, $false, $false, $false, $false, $true, 1, $false, "54×66=3564", 2) | Out-Null
$d.Content.Find.Execute("34×39=1326", $true, $false, $false, $false, $false, $true, 1, $false, "51×88=4488", 2) | Out-Null
$d.Content.Find.Execute("54×12=648", $true, $false, $false, $false, $false, $true, 1, $false, "89×47=4183", 2) | Out-Null
$d.Content.Find.Execute("35×26=910", $true, $false, $false, $false, $false, $true, 1, $false, "68×78=5304", 2) | Out-Null
$d.Content.Find.Execute("55×96=5280", $true, $false, $false, $false, $false, $true, 1, $false, "36×26=936", 2) | Out-Null
$d.Content.Find.Execute("87×34=2958", $true, $false, $false, $false, $false, $true, 1, $false, "56×34=1904", 2) | Out-Null
$d.Content.Find.Execute("21×27=567", $true, $false, $false, $false, $false, $true, 1, $false, "32×44=1408", 2) | Out-Null
$d.Content.Find.Execute("65×83=5395", $true, $false, $false, $false, $false, $true, 1, $false, "77×35=2695", 2) | Out-Null
$d.Content.Find.Execute("38×61=2318", $true, $false, $false, $false, $false, $true, 1, $false, "24×55=1320", 2) | Out-Null
$d.Content.Find.Execute("16×85=1360", $true, $false, $false, $false, $false, $true, 1, $false, "54×32=1728", 2) | Out-Null
$d.Content.Find.Execute("53×34=1802", $true, $false, $false, $false, $false, $true, 1, $false, "92×39=3588", 2) | Out-Null
$d.Content.Find.Execute("17×41=697", $true, $false, $false, $false, $false, $true, 1, $false, "46×85=3910", 2) | Out-Null
$d.Content.Find.Execute("14×25=350", $true, $false, $false, $false, $false, $true, 1, $false, "19×34=646", 2) | Out-Null
$d.Content.Find.Execute("63×81=5103", $true, $false, $false, $false, $false, $true, 1, $false, "39×28=1092", 2) | Out-Null
$d.Content.Find.Execute("32×34=1088", $true, $false, $false, $false, $false, $true, 1, $false, "44×38=1672", 2) | Out-Null
$d.Content.Find.Execute("33×72=2376", $true, $false, $false, $false, $false, $true, 1, $false, "85×24=2040", 2) | Out-Null
$d.Content.Find.Execute("96×84=8064", $true, $false, $false, $false, $false, $true, 1, $false, "91×98=8918", 2) | Out-Null
$d.Content.Find.Execute("58×19=1102", $true, $false, $false, $false, $false, $true, 1, $false, "62×92=5704", 2) | Out-Null
$d.Content.Find.Execute("100×47=4700", $true, $false, $false, $false, $false, $true, 1, $false, "46×58=2668", 2) | Out-Null
$d.Content.Find.Execute("88×48=4224", $true, $false, $false, $false, $false, $true, 1, $false, "19×91=1729", 2) | Out-Null
$d.Content.Find.Execute("54×96=5184", $true, $false, $false, $false, $false, $true, 1, $false, "60×10=600", 2) | Out-Null
$d.Content.Find.Execute("31×12=372", $true, $false, $false, $false, $false, $true, 1, $false, "52×85=4420", 2) | Out-Null
$d.Content.Find.Execute("36×90=3240", $true, $false, $false, $false, $false, $true, 1, $false, "45×53=2385", 2) | Out-Null
$d.Content.Find.Execute("15×20=300", $true, $false, $false, $false, $false, $true, 1, $false, "63×20=1260", 2) | Out-Null
$d.Content.Find.Execute("91×16=1456", $true, $false, $false, $false, $false, $true, 1, $false, "95×54=5130", 2) | Out-Null
$d.Content.Find.Execute("42×70=2940", $true, $false, $false, $false, $false, $true, 1, $false, "77×71=5467", 2) | Out-Null
$d.Content.Find.Execute("95×12=1140", $true, $false, $false, $false, $false, $true, 1, $false, "68×80=5440", 2) | Out-Null
$d.Content.Find.Execute("68×95=6460", $true, $false, $false, $false, $false, $true, 1, $false, "41×58=2378", 2) | Out-Null
$d.Content.Find.Execute("18×99=1782", $true, $false, $false, $false, $false, $true, 1, $false, "16×59=944", 2) | Out-Null
$d.Content.Find.Execute("31×86=2666", $true, $false, $false, $false, $false, $true, 1, $false, "45×41=1845", 2) | Out-Null
$d.Content.Find.Execute("32×50=1600", $true, $false, $false, $false, $false, $true, 1, $false, "83×74=6142", 2) | Out-Null
$d.Content.Find.Execute("40×41=1640", $true, $false, $false, $false, $false, $true, 1, $false, "50×69=3450", 2) | Out-Null
$d.Content.Find.Execute("94×91=8554", $true, $false, $false, $false, $false, $true, 1, $false, "82×25=2050", 2) | Out-Null
$d.Content.Find.Execute("33×66=2178", $true, $false, $false, $false, $false, $true, 1, $false, "16×57=912", 2) | Out-Null
$d.Content.Find.Execute("95×29=2755", $true, $false, $false, $false, $false, $true, 1, $false, "53×22=1166", 2) | Out-Null
$d.Content.Find.Execute("43×83=3569", $true, $false, $false, $false, $false, $true, 1, $false, "80×60=4800", 2) | Out-Null
$d.Content.Find.Execute("48×16=768", $true, $false, $false, $false, $false, $true, 1, $false, "62×87=5394", 2) | Out-Null
$d.Content.Find.Execute("13×54=702", $true, $false, $false, $false, $false, $true, 1, $false, "38×62=2356", 2) | Out-Null
$d.Content.Find.Execute("45×45=2025", $true, $false, $false, $false, $false, $true, 1, $false, "65×87=5655", 2) | Out-Null
$d.Content.Find.Execute("74×31=2294", $true, $false, $false, $false, $false, $true, 1, $false, "73×98=7154", 2) | Out-Null
$d.Content.Find.Execute("95×38=3610", $true, $false, $false, $false, $false, $true, 1, $false, "40×23=920", 2) | Out-Null
$d.Content.Find.Execute("35×93=3255", $true, $false, $false, $false, $false, $true, 1, $false, "83×25=2075", 2) | Out-Null
$d.Content.Find.Execute("82×48=3936", $true, $false, $false, $false, $false, $true, 1, $false, "32×64=2048", 2) | Out-Null
$d.Content.Find.Execute("19×16=304", $true, $false, $false, $false, $false, $true, 1, $false, "10×56=560", 2) | Out-Null
$d.Content.Find.Execute("63×67=4221", $true, $false, $false, $false, $false, $true, 1, $false, "55×22=1210", 2) | Out-Null
$d.Content.Find.Execute("35×22=770", $true, $false, $false, $false, $false, $true, 1, $false, "22×100=2200", 2) | Out-Null
$d.Content.Find.Execute("32×94=3008", $true, $false, $false, $false, $false, $true, 1, $false, "19×74=1406", 2) | Out-Null
$d.Content.Find.Execute("15×95=1425", $true, $false, $false, $false, $false, $true, 1, $false, "43×30=1290", 2) | Out-Null
$d.Content.Find.Execute("19×21=399", $true, $false, $false, $false, $false, $true, 1, $false, "10×29=290", 2) | Out-Null
$d.Content.Find.Execute("52×13=676", $true, $false, $false, $false, $false, $true, 1, $false, "42×31=1302", 2) | Out-Null
$d.Content.Find.Execute("78×65=5070", $true, $false, $false, $false, $false, $true, 1, $false, "25×48=1200", 2) | Out-Null
$d.Content.Find.Execute("15×13=195", $true, $false, $false, $false, $false, $true, 1, $false, "49×88=4312", 2) | Out-Null
$d.Content.Find.Execute("65×28=1820", $true, $false, $false, $false, $false, $true, 1, $false, "88×43=3784", 2) | Out-Null
$d.Content.Find.Execute("29×10=290", $true, $false, $false, $false, $false, $true, 1, $false, "35×17=595", 2) | Out-Null
$d.Content.Find.Execute("27×59=1593", $true, $false, $false, $false, $false, $true, 1, $false, "70×97=6790", 2) | Out-Null
